$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.465.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "  +4.32%  "
$ws.Range("D3").Value = "'1.595.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "  +1.85%  "
$ws.Range("E4").Value2 = "  -0.13%  "
$ws.Range("D5").Value = "'214.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +2.15%  "
$ws.Range("D6").Value = "'0.499"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +1.74%  "
$ws.Range("E7").Value2 = "  -0.15%  "
$ws.Range("D8").Value = "'24.05"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  +8.82%  "
$ws.Range("E9").Value2 = "  +1.08%  "
$ws.Range("D10").Value = "'0.0602"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  +0.82%  "
$ws.Range("D11").Value = "'0.0890"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  +2.29%  "
$ws.Range("D12").Value = "'1.820.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  +1.77%  "
$ws.Range("D13").Value = "'1.592.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  +2.11%  "
$ws.Range("D14").Value = "'3.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  +0.82%  "
$ws.Range("D15").Value = "'0.535"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  +3.21%  "
$ws.Range("D16").Value = "'28.468.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  +4.46%  "
$ws.Range("D17").Value = "'63.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  +2.27%  "
$ws.Range("D18").Value = "'233.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  +7.67%  "
$ws.Range("D19").Value = "'0.0₃0712"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  +1.40%  "
$ws.Range("D20").Value = "'7.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +0.85%  "
$ws.Range("E21").Value2 = "  -0.05%  "
$ws.Range("E22").Value2 = "  -0.18%  "
$ws.Range("D23").Value = "'9.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  +2.51%  "
$ws.Range("E24").Value2 = "  +1.10%  "
$ws.Range("D25").Value = "'152.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  -0.25%  "
$ws.Range("D26").Value = "'15.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  +2.07%  "
$ws.Range("D27").Value = "'6.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -0.03%  "
$ws.Range("E28").Value2 = "  +1.50%  "
$ws.Range("E29").Value2 = "  -0.02%  "
$ws.Range("E30").Value2 = "  +0.73%  "
$ws.Range("D31").Value = "'0.0476"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  +1.13%  "
$ws.Range("E32").Value2 = "  +0.36%  "
$ws.Range("E33").Value2 = "  +0.91%  "
$ws.Range("D34").Value = "'1.424.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  -0.67%  "
$ws.Range("D35").Value = "'1.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  -0.86%  "
$ws.Range("E36").Value2 = "  -3.56%  "
$ws.Range("E37").Value2 = "  -0.10%  "
$ws.Range("E38").Value2 = "  +0.51%  "
$ws.Range("E39").Value2 = "  +8.54%  "
$ws.Range("E40").Value2 = "  +2.10%  "
$ws.Range("D41").Value = "'0.822"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  +1.87%  "
$ws.Range("E42").Value2 = "  -3.08%  "
$ws.Range("E43").Value2 = "  -0.18%  "
$ws.Range("B44").Value2 = "WEMIXToken"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'0.980"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  -1.78%  "
$ws.Range("B45").Value2 = "RenderToken"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'1.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  +6.40%  "
$ws.Range("D46").Value = "'64.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  +0.58%  "
$ws.Range("D47").Value = "'1.733.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  +1.80%  "
$ws.Range("D48").Value = "'87.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +1.75%  "
$ws.Range("E49").Value2 = "  +0.41%  "
$ws.Range("E50").Value2 = "  +5.47%  "
$ws.Range("E51").Value2 = "  -0.30%  "
